# ============================================================================
# Applies the "ADDITIONAL SCRAPING" edit described in the commit:
#   - Splits the single "ODI Batting" sheet workbook into three sheets:
#       1) "Player Info"       (new)   - basic player bio info
#       2) "ODI Batting"       (existing, modified)
#            * header D1: MATCH_CARD_LINK -> MATCH_CODE
#            * column D values: full howstat URL -> bare MatchCode number
#            * a few stray empty B cells (rows 15,16,22) fully cleared
#       3) "ODI Batting Extra" (new)   - extra per-match batting stats
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# 0. Grab a handle on the existing (and currently only) worksheet.
# ----------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item(1)

# ----------------------------------------------------------------------
# 1. Insert "ODI Batting Extra" right after "ODI Batting", and
#    "Player Info" right before it, so the final tab order is:
#    Player Info, ODI Batting, ODI Batting Extra
#
#    NOTE: worksheet handles here track *position*, not identity, so as
#    soon as a new sheet is inserted before/after one we're already
#    holding, the old variable silently starts pointing at whatever now
#    sits in that slot. To stay safe, every sheet is re-fetched **by
#    name** right before it's used for anything past this point.
# ----------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Add($null, $wsBatting)
$wsExtra.Name = "ODI Batting Extra"

$wsBatting = $wb.Worksheets.Item("ODI Batting")
$wsInfo = $wb.Worksheets.Add($wsBatting)
$wsInfo.Name = "Player Info"

# ----------------------------------------------------------------------
# 2. Populate "Player Info"
# ----------------------------------------------------------------------
$wsInfo = $wb.Worksheets.Item("Player Info")

$infoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $infoHeaders.Length; $c++) {
    $cell = $wsInfo.Cells.Item(1, $c)
    $cell.Value = $infoHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

$infoRow = @("4538", "Evin Lewis", "Left Handed", "Does Not Bowl | Unknown")
for ($c = 1; $c -le $infoRow.Length; $c++) {
    $cell = $wsInfo.Cells.Item(2, $c)
    $cell.NumberFormat = "@"
    $cell.Value = $infoRow[$c - 1]
}

$wsInfo.Range("A1").Select()

# ----------------------------------------------------------------------
# 3. Fix up "ODI Batting"
# ----------------------------------------------------------------------
$wsBatting = $wb.Worksheets.Item("ODI Batting")

# 3a. Header rename: MATCH_CARD_LINK -> MATCH_CODE
$wsBatting.Cells.Item(1, 4).Value = "MATCH_CODE"

# 3b. Column D: turn the full howstat scorecard URL into the bare
#     numeric MatchCode that was on the end of the querystring.
$usedRows = $wsBatting.UsedRange.Rows.Count
for ($r = 2; $r -le $usedRows; $r++) {
    $cell = $wsBatting.Cells.Item($r, 4)
    $url = $cell.Value2
    if ($url -ne $null -and $url -ne "") {
        $code = $url -replace ".*MatchCode=", ""
        $cell.NumberFormat = "@"
        $cell.Value = $code
    }
}

# 3c. A handful of rows had a stray, fully-empty INNING_NUMBER (col B)
#     cell that is removed outright in the new version.
$emptyBRows = @(15, 16, 22)
foreach ($r in $emptyBRows) {
    $wsBatting.Cells.Item($r, 2).ClearContents()
}

# ----------------------------------------------------------------------
# 4. Populate "ODI Batting Extra"
# ----------------------------------------------------------------------
$wsExtra = $wb.Worksheets.Item("ODI Batting Extra")

$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    $cell = $wsExtra.Cells.Item(1, $c)
    $cell.Value = $extraHeaders[$c - 1]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# MATCH_CODE, BATTING_POSITION, NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
# BATTING_POSITION is numeric; everything else (bar MAN_OF_MATCH, which is
# never numeric-looking anyway) must stay plain text, matching the source.
$extraRows = @(
    @("4325", "", "", "", "", "NO"),
    @("4333", "", "", "", "", "NO"),
    @("4348", "2", "6", "2", "18.65%", "NO"),
    @("4359", "", "", "", "", "NO"),
    @("4360", "2", "8", "1", "30.95%", "NO"),
    @("4362", "", "", "", "", "NO"),
    @("4377", "2", "0", "0", "3.55%", "NO"),
    @("4378", "2", "6", "1", "21.86%", "NO"),
    @("4379", "2", "0", "0", "0.40%", "NO"),
    @("4387", "1", "5", "0", "10.71%", "NO"),
    @("4388", "1", "3", "0", "6.67%", "NO"),
    @("4391", "2", "13", "2", "53.80%", "NO"),
    @("4394", "2", "1", "0", "2.89%", "NO"),
    @("4397", "2", "6", "5", "51.26%", "YES"),
    @("4449", "1", "4", "2", "27.54%", "NO"),
    @("4450", "", "", "", "", "NO"),
    @("4451", "1", "2", "0", "4.71%", "NO"),
    @("4483", "1", "0", "0", "", "NO"),
    @("4484", "1", "0", "0", "0.52%", "NO"),
    @("4486", "", "", "", "", "NO")
)

for ($i = 0; $i -lt $extraRows.Length; $i++) {
    $r = $i + 2
    $rowData = $extraRows[$i]

    # A: MATCH_CODE (text)
    $cellA = $wsExtra.Cells.Item($r, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $rowData[0]

    # B: BATTING_POSITION (numeric, left blank when unknown - matches
    #    source, which has no cell value at all for those rows)
    $bVal = $rowData[1]
    if ($bVal -ne "") {
        $wsExtra.Cells.Item($r, 2).Value = [double]$bVal
    }

    # C: NUM_4 (text; left blank when source has no value)
    $cVal = $rowData[2]
    if ($cVal -ne "") {
        $cellC = $wsExtra.Cells.Item($r, 3)
        $cellC.NumberFormat = "@"
        $cellC.Value = $cVal
    }

    # D: NUM_6 (text; left blank when source has no value)
    $dVal = $rowData[3]
    if ($dVal -ne "") {
        $cellD = $wsExtra.Cells.Item($r, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $dVal
    }

    # E: PERCENT_RUNS_OF_TOTAL (text; left blank when source has no value)
    $eVal = $rowData[4]
    if ($eVal -ne "") {
        $cellE = $wsExtra.Cells.Item($r, 5)
        $cellE.NumberFormat = "@"
        $cellE.Value = $eVal
    }

    # F: MAN_OF_MATCH (text, always populated)
    $cellF = $wsExtra.Cells.Item($r, 6)
    $cellF.NumberFormat = "@"
    $cellF.Value = $rowData[5]
}

$wsExtra.Range("A1").Select()

# ----------------------------------------------------------------------
# 5. Leave the view on the first sheet, matching activeTab=0 in the diff.
# ----------------------------------------------------------------------
$wb.Worksheets.Item("Player Info").Activate()
